$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1. Summary paragraph: "Data Scientist" -> "Data scientist", plus two
#    added commas ("cleaning," and "algebraic,").
# ----------------------------------------------------------------------
$d.Content.Find.Execute(
    "Data Scientist familiar with gathering, cleaning and organizing data for use by technical and non-technical personnel. Advanced understanding of statistical, algebraic and other analytical techniques. Highly organized, motivated and diligent with significant background in predictive analytics",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Data scientist familiar with gathering, cleaning, and organizing data for use by technical and non-technical personnel. Advanced understanding of statistical, algebraic, and other analytical techniques. Highly organized, motivated and diligent with significant background in predictive analytics",
    2) | Out-Null

# ----------------------------------------------------------------------
# 2. "Ekohealth, Oakland, CA" -- merge runs that used to be split around
#    a proofing (spell-check) mark; text itself is unchanged.
# ----------------------------------------------------------------------
$d.Content.Find.Execute(
    "Ekohealth, Oakland, CA",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Ekohealth, Oakland, CA",
    2) | Out-Null

# ----------------------------------------------------------------------
# 3. "Helped build the product Eko-core, ..." -- merge runs.
# ----------------------------------------------------------------------
$d.Content.Find.Execute(
    "Helped build the product Eko-core, an FDA-cleared digital stethoscope attachment device, saving monthly cost for patients with arteriovenous fistula (AVF)",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Helped build the product Eko-core, an FDA-cleared digital stethoscope attachment device, saving monthly cost for patients with arteriovenous fistula (AVF)",
    2) | Out-Null

# ----------------------------------------------------------------------
# 4. "Productionalized customer-facing ..." -- merge runs.
# ----------------------------------------------------------------------
$d.Content.Find.Execute(
    "Productionalized customer-facing python-based analysis pipeline using AWS cloud services",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Productionalized customer-facing python-based analysis pipeline using AWS cloud services",
    2) | Out-Null

# ----------------------------------------------------------------------
# 5. "... using dplyr and tidyverse" -- merge runs.
# ----------------------------------------------------------------------
$d.Content.Find.Execute(
    "Migrated data from SAS to R and performed EDA using dplyr and tidyverse",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Migrated data from SAS to R and performed EDA using dplyr and tidyverse",
    2) | Out-Null

# ----------------------------------------------------------------------
# 6. The "Decreased the data dimensionality ..." bullet becomes a CNN /
#    Keras bullet, and a brand-new bullet (new paragraph) is added right
#    after it carrying a reworded "Decreased the data dimensionality ..."
#    line (now mentioning the 23% improvement instead of the old MSE
#    number).
# ----------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Decreased the data dimensionality using principal component analysis*") {
        # Insert a new (empty) bullet paragraph right after this one --
        # it inherits this paragraph's pPr (style / numbering / tabs).
        $p.Range.InsertParagraphAfter()
        $newPara = $p.Next()
        $newRange = $newPara.Range
        $newRange.End = $newRange.End - 1
        $newRange.Text = "Decreased the data dimensionality using principal component analysis (PCA) and improved prediction by 23% training a generalized linear model (GLM)"

        # Now overwrite the original bullet with the new CNN / Keras text.
        $origRange = $p.Range
        $origRange.End = $origRange.End - 1
        $origRange.Text = "Trained convolutional neural networks (CNN) using Keras for BMI prediction with 25.45 mean squared error (MSE)"
        break
    }
}

# ----------------------------------------------------------------------
# 7. "Hosted R Shiny website ... using ggplot2 and plotly" -- merge runs
#    around the trailing "plotly" word.
# ----------------------------------------------------------------------
$d.Content.Find.Execute(
    "comparing machine learning algorithms (PCA, k-means, UMAP, and t-SNE) & visualized clustering results using ggplot2 and plotly",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "comparing machine learning algorithms (PCA, k-means, UMAP, and t-SNE) & visualized clustering results using ggplot2 and plotly",
    2) | Out-Null

# ----------------------------------------------------------------------
# 8. "Processed Amazon Food Review data using pandas, NumPy, and dfply
#    in Python" -- merge runs.
# ----------------------------------------------------------------------
$d.Content.Find.Execute(
    "Processed Amazon Food Review data using pandas, NumPy, and dfply in Python",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Processed Amazon Food Review data using pandas, NumPy, and dfply in Python",
    2) | Out-Null

# ----------------------------------------------------------------------
# 9. "Python, R (RShiny), SQL, Shell scripting" -- merge runs around
#    "Shiny".
# ----------------------------------------------------------------------
$d.Content.Find.Execute(
    "Python, R (RShiny), SQL, Shell scripting",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Python, R (RShiny), SQL, Shell scripting",
    2) | Out-Null

# ----------------------------------------------------------------------
# 10. "Tableau, Matplotlib, Seaborn, ggplot2, plotly" -- merge runs.
# ----------------------------------------------------------------------
$d.Content.Find.Execute(
    "Tableau, Matplotlib, Seaborn, ggplot2, plotly",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Tableau, Matplotlib, Seaborn, ggplot2, plotly",
    2) | Out-Null

# ----------------------------------------------------------------------
# 11. "Pandas, NumPy, SciPy, NLTK, scikit-learn, Tidyverse" -- merge
#     runs.
# ----------------------------------------------------------------------
$d.Content.Find.Execute(
    "Pandas, NumPy, SciPy, NLTK, scikit-learn, Tidyverse",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Pandas, NumPy, SciPy, NLTK, scikit-learn, Tidyverse",
    2) | Out-Null

Write-Output "done"
